$d = $word.ActiveDocument

$d.Content.Find.Execute("65-58=", $true, $false, $false, $false, $false, $true, 1, $false, "5+63=", 2) | Out-Null
$d.Content.Find.Execute("32+59=", $true, $false, $false, $false, $false, $true, 1, $false, "16+57=", 2) | Out-Null
$d.Content.Find.Execute("5+36=", $true, $false, $false, $false, $false, $true, 1, $false, "68-45=", 2) | Out-Null
$d.Content.Find.Execute("77-60=", $true, $false, $false, $false, $false, $true, 1, $false, "92+5=", 2) | Out-Null
$d.Content.Find.Execute("72-59=", $true, $false, $false, $false, $false, $true, 1, $false, "83-39=", 2) | Out-Null
$d.Content.Find.Execute("95-9=", $true, $false, $false, $false, $false, $true, 1, $false, "83-79=", 2) | Out-Null
$d.Content.Find.Execute("30+62=", $true, $false, $false, $false, $false, $true, 1, $false, "11+88=", 2) | Out-Null
$d.Content.Find.Execute("6+76=", $true, $false, $false, $false, $false, $true, 1, $false, "72-37=", 2) | Out-Null
$d.Content.Find.Execute("55+35=", $true, $false, $false, $false, $false, $true, 1, $false, "22+46=", 2) | Out-Null
$d.Content.Find.Execute("5+0=", $true, $false, $false, $false, $false, $true, 1, $false, "36-11=", 2) | Out-Null
$d.Content.Find.Execute("39+38=", $true, $false, $false, $false, $false, $true, 1, $false, "48+46=", 2) | Out-Null
$d.Content.Find.Execute("42+37=", $true, $false, $false, $false, $false, $true, 1, $false, "88-28=", 2) | Out-Null
$d.Content.Find.Execute("89-84=", $true, $false, $false, $false, $false, $true, 1, $false, "2+2=", 2) | Out-Null
$d.Content.Find.Execute("22+76=", $true, $false, $false, $false, $false, $true, 1, $false, "38+50=", 2) | Out-Null
$d.Content.Find.Execute("2+84=", $true, $false, $false, $false, $false, $true, 1, $false, "47-8=", 2) | Out-Null
$d.Content.Find.Execute("78-5=", $true, $false, $false, $false, $false, $true, 1, $false, "56-7=", 2) | Out-Null
$d.Content.Find.Execute("91-28=", $true, $false, $false, $false, $false, $true, 1, $false, "55-25=", 2) | Out-Null
$d.Content.Find.Execute("0+13=", $true, $false, $false, $false, $false, $true, 1, $false, "28+32=", 2) | Out-Null
$d.Content.Find.Execute("33+14=", $true, $false, $false, $false, $false, $true, 1, $false, "6+5=", 2) | Out-Null
$d.Content.Find.Execute("59+39=", $true, $false, $false, $false, $false, $true, 1, $false, "79-37=", 2) | Out-Null
$d.Content.Find.Execute("7-2=", $true, $false, $false, $false, $false, $true, 1, $false, "7+63=", 2) | Out-Null
$d.Content.Find.Execute("91-1=", $true, $false, $false, $false, $false, $true, 1, $false, "79+15=", 2) | Out-Null
$d.Content.Find.Execute("73-31=", $true, $false, $false, $false, $false, $true, 1, $false, "73-28=", 2) | Out-Null
$d.Content.Find.Execute("11+84=", $true, $false, $false, $false, $false, $true, 1, $false, "66+15=", 2) | Out-Null
$d.Content.Find.Execute("95+3=", $true, $false, $false, $false, $false, $true, 1, $false, "53-30=", 2) | Out-Null
$d.Content.Find.Execute("77-57=", $true, $false, $false, $false, $false, $true, 1, $false, "80-64=", 2) | Out-Null
$d.Content.Find.Execute("42+14=", $true, $false, $false, $false, $false, $true, 1, $false, "62-33=", 2) | Out-Null
$d.Content.Find.Execute("22+35=", $true, $false, $false, $false, $false, $true, 1, $false, "35-30=", 2) | Out-Null
$d.Content.Find.Execute("38-6=", $true, $false, $false, $false, $false, $true, 1, $false, "18-3=", 2) | Out-Null
$d.Content.Find.Execute("90-10=", $true, $false, $false, $false, $false, $true, 1, $false, "56-31=", 2) | Out-Null
$d.Content.Find.Execute("77+6=", $true, $false, $false, $false, $false, $true, 1, $false, "96-14=", 2) | Out-Null
$d.Content.Find.Execute("78-4=", $true, $false, $false, $false, $false, $true, 1, $false, "21+1=", 2) | Out-Null
$d.Content.Find.Execute("66-25=", $true, $false, $false, $false, $false, $true, 1, $false, "64+28=", 2) | Out-Null
$d.Content.Find.Execute("23+70=", $true, $false, $false, $false, $false, $true, 1, $false, "57-7=", 2) | Out-Null
$d.Content.Find.Execute("66-33=", $true, $false, $false, $false, $false, $true, 1, $false, "86-38=", 2) | Out-Null
$d.Content.Find.Execute("59-44=", $true, $false, $false, $false, $false, $true, 1, $false, "9+73=", 2) | Out-Null
$d.Content.Find.Execute("69-18=", $true, $false, $false, $false, $false, $true, 1, $false, "22-11=", 2) | Out-Null
$d.Content.Find.Execute("66+18=", $true, $false, $false, $false, $false, $true, 1, $false, "45+9=", 2) | Out-Null
$d.Content.Find.Execute("50-31=", $true, $false, $false, $false, $false, $true, 1, $false, "72-48=", 2) | Out-Null
$d.Content.Find.Execute("44+36=", $true, $false, $false, $false, $false, $true, 1, $false, "10+88=", 2) | Out-Null
$d.Content.Find.Execute("88-68=", $true, $false, $false, $false, $false, $true, 1, $false, "54-4=", 2) | Out-Null
$d.Content.Find.Execute("58+11=", $true, $false, $false, $false, $false, $true, 1, $false, "73-53=", 2) | Out-Null
$d.Content.Find.Execute("44-37=", $true, $false, $false, $false, $false, $true, 1, $false, "35+56=", 2) | Out-Null
$d.Content.Find.Execute("24+1=", $true, $false, $false, $false, $false, $true, 1, $false, "15+34=", 2) | Out-Null
$d.Content.Find.Execute("92-8=", $true, $false, $false, $false, $false, $true, 1, $false, "75-20=", 2) | Out-Null
$d.Content.Find.Execute("4+65=", $true, $false, $false, $false, $false, $true, 1, $false, "83-16=", 2) | Out-Null
$d.Content.Find.Execute("6+82=", $true, $false, $false, $false, $false, $true, 1, $false, "58+25=", 2) | Out-Null
$d.Content.Find.Execute("49-31=", $true, $false, $false, $false, $false, $true, 1, $false, "0+26=", 2) | Out-Null
$d.Content.Find.Execute("14+2=", $true, $false, $false, $false, $false, $true, 1, $false, "85-32=", 2) | Out-Null
$d.Content.Find.Execute("36-6=", $true, $false, $false, $false, $false, $true, 1, $false, "81-74=", 2) | Out-Null
$d.Content.Find.Execute("82-73=", $true, $false, $false, $false, $false, $true, 1, $false, "26-12=", 2) | Out-Null
$d.Content.Find.Execute("0+93=", $true, $false, $false, $false, $false, $true, 1, $false, "36+12=", 2) | Out-Null
$d.Content.Find.Execute("3+26=", $true, $false, $false, $false, $false, $true, 1, $false, "90+8=", 2) | Out-Null
$d.Content.Find.Execute("73-30=", $true, $false, $false, $false, $false, $true, 1, $false, "73-45=", 2) | Out-Null
$d.Content.Find.Execute("74-57=", $true, $false, $false, $false, $false, $true, 1, $false, "18+44=", 2) | Out-Null
$d.Content.Find.Execute("21+32=", $true, $false, $false, $false, $false, $true, 1, $false, "90-60=", 2) | Out-Null
$d.Content.Find.Execute("92-72=", $true, $false, $false, $false, $false, $true, 1, $false, "30-19=", 2) | Out-Null
$d.Content.Find.Execute("37-7=", $true, $false, $false, $false, $false, $true, 1, $false, "45-16=", 2) | Out-Null
$d.Content.Find.Execute("20+7=", $true, $false, $false, $false, $false, $true, 1, $false, "15+26=", 2) | Out-Null
$d.Content.Find.Execute("33-13=", $true, $false, $false, $false, $false, $true, 1, $false, "86-40=", 2) | Out-Null
$d.Content.Find.Execute("93-72=", $true, $false, $false, $false, $false, $true, 1, $false, "86+0=", 2) | Out-Null
$d.Content.Find.Execute("80-41=", $true, $false, $false, $false, $false, $true, 1, $false, "42+30=", 2) | Out-Null
$d.Content.Find.Execute("84-71=", $true, $false, $false, $false, $false, $true, 1, $false, "3+1=", 2) | Out-Null
$d.Content.Find.Execute("86-54=", $true, $false, $false, $false, $false, $true, 1, $false, "22+14=", 2) | Out-Null
$d.Content.Find.Execute("98-76=", $true, $false, $false, $false, $false, $true, 1, $false, "99-68=", 2) | Out-Null
$d.Content.Find.Execute("13+5=", $true, $false, $false, $false, $false, $true, 1, $false, "89-36=", 2) | Out-Null
$d.Content.Find.Execute("49-39=", $true, $false, $false, $false, $false, $true, 1, $false, "12+54=", 2) | Out-Null
$d.Content.Find.Execute("69-23=", $true, $false, $false, $false, $false, $true, 1, $false, "9-7=", 2) | Out-Null
$d.Content.Find.Execute("24+45=", $true, $false, $false, $false, $false, $true, 1, $false, "38+26=", 2) | Out-Null
$d.Content.Find.Execute("85-50=", $true, $false, $false, $false, $false, $true, 1, $false, "63-22=", 2) | Out-Null
$d.Content.Find.Execute("45+18=", $true, $false, $false, $false, $false, $true, 1, $false, "50-4=", 2) | Out-Null
$d.Content.Find.Execute("57+18=", $true, $false, $false, $false, $false, $true, 1, $false, "79-44=", 2) | Out-Null
$d.Content.Find.Execute("34+51=", $true, $false, $false, $false, $false, $true, 1, $false, "37+61=", 2) | Out-Null
$d.Content.Find.Execute("35-10=", $true, $false, $false, $false, $false, $true, 1, $false, "32+58=", 2) | Out-Null
$d.Content.Find.Execute("35+25=", $true, $false, $false, $false, $false, $true, 1, $false, "32+26=", 2) | Out-Null
$d.Content.Find.Execute("90-55=", $true, $false, $false, $false, $false, $true, 1, $false, "26+35=", 2) | Out-Null
$d.Content.Find.Execute("92-85=", $true, $false, $false, $false, $false, $true, 1, $false, "67+31=", 2) | Out-Null
$d.Content.Find.Execute("86-8=", $true, $false, $false, $false, $false, $true, 1, $false, "29-10=", 2) | Out-Null
$d.Content.Find.Execute("75-21=", $true, $false, $false, $false, $false, $true, 1, $false, "66-6=", 2) | Out-Null
$d.Content.Find.Execute("28+47=", $true, $false, $false, $false, $false, $true, 1, $false, "96-82=", 2) | Out-Null
$d.Content.Find.Execute("90-34=", $true, $false, $false, $false, $false, $true, 1, $false, "34-14=", 2) | Out-Null
$d.Content.Find.Execute("58-33=", $true, $false, $false, $false, $false, $true, 1, $false, "80-60=", 2) | Out-Null
$d.Content.Find.Execute("71-64=", $true, $false, $false, $false, $false, $true, 1, $false, "54+31=", 2) | Out-Null
$d.Content.Find.Execute("43-10=", $true, $false, $false, $false, $false, $true, 1, $false, "23+44=", 2) | Out-Null
$d.Content.Find.Execute("31-21=", $true, $false, $false, $false, $false, $true, 1, $false, "79+15=", 2) | Out-Null
$d.Content.Find.Execute("84-23=", $true, $false, $false, $false, $false, $true, 1, $false, "75+1=", 2) | Out-Null
$d.Content.Find.Execute("41-39=", $true, $false, $false, $false, $false, $true, 1, $false, "4+2=", 2) | Out-Null
$d.Content.Find.Execute("77-62=", $true, $false, $false, $false, $false, $true, 1, $false, "99-61=", 2) | Out-Null
$d.Content.Find.Execute("85+5=", $true, $false, $false, $false, $false, $true, 1, $false, "86-31=", 2) | Out-Null
$d.Content.Find.Execute("72-3=", $true, $false, $false, $false, $false, $true, 1, $false, "4+3=", 2) | Out-Null
$d.Content.Find.Execute("9+67=", $true, $false, $false, $false, $false, $true, 1, $false, "37-14=", 2) | Out-Null
$d.Content.Find.Execute("5+7=", $true, $false, $false, $false, $false, $true, 1, $false, "70-12=", 2) | Out-Null
$d.Content.Find.Execute("54-8=", $true, $false, $false, $false, $false, $true, 1, $false, "67+1=", 2) | Out-Null
$d.Content.Find.Execute("97-0=", $true, $false, $false, $false, $false, $true, 1, $false, "79-45=", 2) | Out-Null
$d.Content.Find.Execute("7+54=", $true, $false, $false, $false, $false, $true, 1, $false, "50-49=", 2) | Out-Null
$d.Content.Find.Execute("71-27=", $true, $false, $false, $false, $false, $true, 1, $false, "54+18=", 2) | Out-Null
$d.Content.Find.Execute("57-10=", $true, $false, $false, $false, $false, $true, 1, $false, "53-34=", 2) | Out-Null
$d.Content.Find.Execute("97-3=", $true, $false, $false, $false, $false, $true, 1, $false, "89-59=", 2) | Out-Null
$d.Content.Find.Execute("59+29=", $true, $false, $false, $false, $false, $true, 1, $false, "28+38=", 2) | Out-Null
$d.Content.Find.Execute("64+6=", $true, $false, $false, $false, $false, $true, 1, $false, "49-46=", 2) | Out-Null
